{"js": "// Insert four new bullet paragraphs describing the Python boundary\n// estimation algorithm work right after the \"Full-Stack Development and\n// Data Engineering\" sub-heading paragraph under the Siege Analytics entry,\n// and before the existing \"Lead comprehensive research initiatives...\"\n// bullet.\n\nconst body = context.document.body;\n\n// Locate the anchor paragraph by its unique text.\nconst results = body.search(\"Full-Stack Development and Data Engineering\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\n    \"Could not find the 'Full-Stack Development and Data Engineering' paragraph\"\n  );\n}\n\nlet anchorParagraph = results.items[0].paragraphs.getFirst();\n\nconst newBullets = [\n  \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\u2022 Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times\",\n  \"\u2022 Architected systems supporting 2,500+ concurrent users conducting redistricting analysis\",\n  \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\n// Insert each bullet right after the previous one so they end up in order\n// immediately following the anchor paragraph.\nfor (const bulletText of newBullets) {\n  anchorParagraph = anchorParagraph.insertParagraph(bulletText, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert four new bullet paragraphs describing the Python boundary\n# estimation algorithm work right after the \"Full-Stack Development and\n# Data Engineering\" sub-heading paragraph under the Siege Analytics entry,\n# and before the existing \"Lead comprehensive research initiatives...\"\n# bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its exact text (paragraph Range.Text\n# includes the trailing paragraph mark, chr(13), so trim it before\n# comparing).\n$anchorIndex = 0\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $paraText = $p.Range.Text.TrimEnd([char]13)\n    if ($paraText -eq \"Full-Stack Development and Data Engineering\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq 0) {\n    Write-Output \"Could not find the 'Full-Stack Development and Data Engineering' paragraph\"\n} else {\n    $newBullets = @(\n        \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n        \"\u2022 Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times\",\n        \"\u2022 Architected systems supporting 2,500+ concurrent users conducting redistricting analysis\",\n        \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n    )\n\n    $currentIndex = $anchorIndex\n    foreach ($bulletText in $newBullets) {\n        $currentParagraph = $d.Paragraphs.Item($currentIndex)\n        $currentParagraph.Range.InsertParagraphAfter()\n        $currentIndex = $currentIndex + 1\n        $newParagraph = $d.Paragraphs.Item($currentIndex)\n        $newParagraph.Range.Text = $bulletText\n    }\n}\n"}
